$d = $word.ActiveDocument

$pairs = @(
    @("510×7=3570", "636×8=5088"),
    @("102×8=816", "107×8=856"),
    @("757×8=6056", "430×6=2580"),
    @("505×4=2020", "370×2=740"),
    @("112×9=1008", "296×4=1184"),
    @("602×7=4214", "500×3=1500"),
    @("390×2=780", "368×7=2576"),
    @("998×3=2994", "484×4=1936"),
    @("336×8=2688", "385×2=770"),
    @("199×5=995", "374×9=3366"),
    @("467×6=2802", "475×9=4275"),
    @("428×3=1284", "476×5=2380"),
    @("270×3=810", "193×6=1158"),
    @("512×4=2048", "564×8=4512"),
    @("787×2=1574", "248×3=744"),
    @("799×5=3995", "738×7=5166"),
    @("271×8=2168", "315×2=630"),
    @("522×4=2088", "114×6=684"),
    @("685×4=2740", "147×6=882"),
    @("166×4=664", "393×3=1179"),
    @("739×7=5173", "822×5=4110"),
    @("509×8=4072", "807×3=2421"),
    @("562×8=4496", "691×9=6219"),
    @("495×5=2475", "611×3=1833"),
    @("364×6=2184", "222×5=1110")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
